$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-21 13:08:42"

$wsZhCn.Range("H2").Value = "2016-08-21 13:08:38"
$wsZhCn.Range("K2").Value = "2016-08-21 13:08:56"

$wsDeDe.Range("H2").Value = "2016-08-21 13:08:42"
$wsDeDe.Range("K2").Value = "2016-08-21 13:09:07"
